$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: insert a new column I (becomes "Elapsed Time") ---
$ws.Columns("I:I").Insert()

# --- Step 2: insert a new column M (first of the two new Fairness columns) ---
$ws.Columns("M:M").Insert()

# --- Step 3: insert a new column O (second new Fairness column, after the old Fairness col which is now N) ---
$ws.Columns("O:O").Insert()

# --- Headers (order matches the shared-string table in the target file) ---
$ws.Range("M1").Value = "Fairness(QWT)"
$ws.Range("N1").Value = "Fairness(RT)"
$ws.Range("O1").Value = "Fairness(ET)"
$ws.Range("I1").Value = "Elapsed Time"

# --- New data formulas for column I (Elapsed Time), mirroring H but referencing C ---
$ws.Range("I2").Formula = "=(C2-`$B`$12)/100"
$ws.Range("I3:I11").Formula = "=(C3-`$B`$12)/100"

# --- Re-assert the (unchanged) formulas in the columns that shifted during the
#     inserts above, so the shared-formula grouping is rebuilt across the full
#     column range exactly like it was before the insert ---
$ws.Range("L3:L11").Formula = "=((E3+F3)/`$J`$12)*100"
$ws.Range("P3:P11").Formula = "=F3"

# --- New summary formulas ---
$ws.Range("M12").Formula = "=_xlfn.STDEV.P(G2:G11)"
$ws.Range("O12").Formula = "=_xlfn.STDEV.P(I2:I11)"

# --- Selection matches the saved file's cursor position ---
$null = $ws.Range("I16").Select()
